$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.47086066666667
$ws.Range("H2").Value = 31.412582
$ws.Range("I2").Value = 0.340259118787888
$ws.Range("J2").Value = 0.340259118787888
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.709791333333333
$ws.Range("N2").Value = 11.129374
$ws.Range("O2").Value = 0.4283284425582907
$ws.Range("P2").Value = 0.4283284425582907
$ws.Range("Q2").Value = 38.84470815374089
$ws.Range("R2").Value = 349.602373383668
$ws.Range("S2").Value = 0.1457426584166725
$ws.Range("T2").Value = 0.1457426584166725
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.47086066666667
$ws.Range("H3").Value = 31.412582
$ws.Range("I3").Value = 0.340259118787888
$ws.Range("J3").Value = 0.340259118787888
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.599001333333334
$ws.Range("N3").Value = 10.797004
$ws.Range("O3").Value = 0.4155367505499981
$ws.Range("P3").Value = 0.4155367505499982
$ws.Range("Q3").Value = 37.68464150048089
$ws.Range("R3").Value = 339.161773504328
$ws.Range("S3").Value = 0.1413901685661248
$ws.Range("T3").Value = 0.1413901685661248
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.47086066666667
$ws.Range("H4").Value = 31.412582
$ws.Range("I4").Value = 0.340259118787888
$ws.Range("J4").Value = 0.340259118787888
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.352297666666667
$ws.Range("N4").Value = 4.056893
$ws.Range("O4").Value = 0.1561348068917112
$ws.Range("P4").Value = 0.1561348068917112
$ws.Range("Q4").Value = 14.15972044752511
$ws.Range("R4").Value = 127.437484027726
$ws.Range("S4").Value = 0.05312629180509072
$ws.Range("T4").Value = 0.05312629180509072
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.84313566666667
$ws.Range("H5").Value = 38.529407
$ws.Range("I5").Value = 0.417348120993043
$ws.Range("J5").Value = 0.4173481209930429
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.709791333333333
$ws.Range("N5").Value = 11.129374
$ws.Range("O5").Value = 0.4283284425582907
$ws.Range("P5").Value = 0.4283284425582907
$ws.Range("Q5").Value = 47.64535338902422
$ws.Range("R5").Value = 428.808180501218
$ws.Range("S5").Value = 0.1787620706695791
$ws.Range("T5").Value = 0.1787620706695791
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.84313566666667
$ws.Range("H6").Value = 38.529407
$ws.Range("I6").Value = 0.417348120993043
$ws.Range("J6").Value = 0.4173481209930429
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.599001333333334
$ws.Range("N6").Value = 10.797004
$ws.Range("O6").Value = 0.4155367505499981
$ws.Range("P6").Value = 0.4155367505499982
$ws.Range("Q6").Value = 46.22246238851423
$ws.Range("R6").Value = 416.002161496628
$ws.Range("S6").Value = 0.1734234820455965
$ws.Range("T6").Value = 0.1734234820455965
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.84313566666667
$ws.Range("H7").Value = 38.529407
$ws.Range("I7").Value = 0.417348120993043
$ws.Range("J7").Value = 0.4173481209930429
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.352297666666667
$ws.Range("N7").Value = 4.056893
$ws.Range("O7").Value = 0.1561348068917112
$ws.Range("P7").Value = 0.1561348068917112
$ws.Range("Q7").Value = 17.36774239471678
$ws.Range("R7").Value = 156.309681552451
$ws.Range("S7").Value = 0.06516256827786729
$ws.Range("T7").Value = 0.06516256827786727
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.820871333333334
$ws.Range("H8").Value = 5.462614
$ws.Range("I8").Value = 0.0591706923651924
$ws.Range("J8").Value = 0.05917069236519239
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.709791333333333
$ws.Range("N8").Value = 11.129374
$ws.Range("O8").Value = 0.4283284425582907
$ws.Range("P8").Value = 0.4283284425582907
$ws.Range("Q8").Value = 6.755052691515112
$ws.Range("R8").Value = 60.79547422363601
$ws.Range("S8").Value = 0.0253444905058786
$ws.Range("T8").Value = 0.0253444905058786
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.820871333333334
$ws.Range("H9").Value = 5.462614
$ws.Range("I9").Value = 0.0591706923651924
$ws.Range("J9").Value = 0.05917069236519239
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.599001333333334
$ws.Range("N9").Value = 10.797004
$ws.Range("O9").Value = 0.4155367505499981
$ws.Range("P9").Value = 0.4155367505499982
$ws.Range("Q9").Value = 6.553318356495113
$ws.Range("R9").Value = 58.97986520845601
$ws.Range("S9").Value = 0.02458759723322563
$ws.Range("T9").Value = 0.02458759723322563
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.820871333333334
$ws.Range("H10").Value = 5.462614
$ws.Range("I10").Value = 0.0591706923651924
$ws.Range("J10").Value = 0.05917069236519239
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.352297666666667
$ws.Range("N10").Value = 4.056893
$ws.Range("O10").Value = 0.1561348068917112
$ws.Range("P10").Value = 0.1561348068917112
$ws.Range("Q10").Value = 2.462360055366889
$ws.Range("R10").Value = 22.161240498302
$ws.Range("S10").Value = 0.009238604626088165
$ws.Range("T10").Value = 0.009238604626088165
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.169900666666666
$ws.Range("H11").Value = 6.509701999999999
$ws.Range("I11").Value = 0.07051268393320077
$ws.Range("J11").Value = 0.07051268393320077
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.709791333333333
$ws.Range("N11").Value = 11.129374
$ws.Range("O11").Value = 0.4283284425582907
$ws.Range("P11").Value = 0.4283284425582907
$ws.Range("Q11").Value = 8.049878687394221
$ws.Range("R11").Value = 72.448908186548
$ws.Range("S11").Value = 0.03020258808971289
$ws.Range("T11").Value = 0.0302025880897129
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.169900666666666
$ws.Range("H12").Value = 6.509701999999999
$ws.Range("I12").Value = 0.07051268393320077
$ws.Range("J12").Value = 0.07051268393320077
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.599001333333334
$ws.Range("N12").Value = 10.797004
$ws.Range("O12").Value = 0.4155367505499981
$ws.Range("P12").Value = 0.4155367505499982
$ws.Range("Q12").Value = 7.809475392534222
$ws.Range("R12").Value = 70.28527853280799
$ws.Range("S12").Value = 0.02930061155416131
$ws.Range("T12").Value = 0.02930061155416131
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.169900666666666
$ws.Range("H13").Value = 6.509701999999999
$ws.Range("I13").Value = 0.07051268393320077
$ws.Range("J13").Value = 0.07051268393320077
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.352297666666667
$ws.Range("N13").Value = 4.056893
$ws.Range("O13").Value = 0.1561348068917112
$ws.Range("P13").Value = 0.1561348068917112
$ws.Range("Q13").Value = 2.934351608431777
$ws.Range("R13").Value = 26.40916447588599
$ws.Range("S13").Value = 0.01100948428932657
$ws.Range("T13").Value = 0.01100948428932657
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.468428
$ws.Range("H14").Value = 10.405284
$ws.Range("I14").Value = 0.1127093839206758
$ws.Range("J14").Value = 0.1127093839206758
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.709791333333333
$ws.Range("N14").Value = 11.129374
$ws.Range("O14").Value = 0.4283284425582907
$ws.Range("P14").Value = 0.4283284425582907
$ws.Range("Q14").Value = 12.86714413469067
$ws.Range("R14").Value = 115.804297212216
$ws.Range("S14").Value = 0.04827663487644752
$ws.Range("T14").Value = 0.04827663487644752
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.468428
$ws.Range("H15").Value = 10.405284
$ws.Range("I15").Value = 0.1127093839206758
$ws.Range("J15").Value = 0.1127093839206758
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.599001333333334
$ws.Range("N15").Value = 10.797004
$ws.Range("O15").Value = 0.4155367505499981
$ws.Range("P15").Value = 0.4155367505499982
$ws.Range("Q15").Value = 12.48287699657067
$ws.Range("R15").Value = 112.345892969136
$ws.Range("S15").Value = 0.04683489115088983
$ws.Range("T15").Value = 0.04683489115088984
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.468428
$ws.Range("H16").Value = 10.405284
$ws.Range("I16").Value = 0.1127093839206758
$ws.Range("J16").Value = 0.1127093839206758
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.352297666666667
$ws.Range("N16").Value = 4.056893
$ws.Range("O16").Value = 0.1561348068917112
$ws.Range("P16").Value = 0.1561348068917112
$ws.Range("Q16").Value = 4.690347091401333
$ws.Range("R16").Value = 42.213123822612
$ws.Range("S16").Value = 0.06516256827786729
$ws.Range("T16").Value = 0.06516256827786727
